# Edit script:
#  1. Remove stray empty inline-string cells B2:B4 on "ODI Batting".
#  2. Add a new worksheet "ODI Batting Extra" (4th sheet, after "ODI Bowling")
#     with per-match batting extras data.

$wb = $excel.ActiveWorkbook

# --- 1. Clear the (empty) B2:B4 cells on "ODI Batting" ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2:B4").ClearContents()

# --- 2. Add new worksheet "ODI Batting Extra" as the last (4th) sheet ---
# Copy an existing data sheet so the new sheet inherits the same sheetPr /
# pageMargins / header-style boilerplate used throughout this workbook, then
# wipe it down to a blank 6-column sheet before filling in the real data.
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Copy($null, $bowling)
$extra = $wb.Worksheets.Item($wb.Worksheets.Count)
$extra.Name = "ODI Batting Extra"

# Drop the copied data/formatting outside the 6 columns we need, and clear
# the old data rows entirely (content + formatting) - row 1's formatting
# (bold header style) is kept and simply overwritten with new header text.
$extra.Range("A2:G9").Clear()
$extra.Columns.Item(7).Clear()

# Header row
$extra.Range("A1").Value = "MATCH_CODE"
$extra.Range("B1").Value = "BATTING_POSITION"
$extra.Range("C1").Value = "NUM_4"
$extra.Range("D1").Value = "NUM_6"
$extra.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$extra.Range("F1").Value = "MAN_OF_MATCH"

# Data rows
# MATCH_CODE (A) is stored as text; prefix numeric-looking strings with an
# apostrophe so Excel keeps them as text instead of coercing to numbers.
# BATTING_POSITION (B) is a genuine number where present.
$extra.Range("A2").Value = "'4388"
$extra.Range("B2").Value = 10
$extra.Range("F2").Value = "NO"

$extra.Range("A3").Value = "'4399"
$extra.Range("B3").Value = 10
$extra.Range("F3").Value = "NO"

$extra.Range("A4").Value = "'4400"
$extra.Range("F4").Value = "NO"

$extra.Range("A5").Value = "'4406"
$extra.Range("B5").Value = 9
$extra.Range("C5").Value = "'2"
$extra.Range("D5").Value = "'0"
$extra.Range("E5").Value = "'2.70%"
$extra.Range("F5").Value = "NO"

$extra.Range("A6").Value = "'4410"
$extra.Range("F6").Value = "NO"

$extra.Range("A7").Value = "'4435"
$extra.Range("F7").Value = "NO"

$extra.Range("A8").Value = "'4436"
$extra.Range("B8").Value = 10
$extra.Range("C8").Value = "'0"
$extra.Range("D8").Value = "'0"
$extra.Range("E8").Value = "'6.67%"
$extra.Range("F8").Value = "NO"

$extra.Range("A9").Value = "'4485"

# Restore the originally active sheet/tab (creating/copying sheets above
# shifts focus onto the new sheet as a side effect).
$wb.Worksheets.Item(1).Select()

